# HIKER-M Update Attendance List[TV]
# Add a new weekly attendance column (J) dated 2021-05-19, mark everyone
# present ("checkmark") by default, and note that Jaime Nunez Delgado left
# early on 2021-05-12 (column I, row 9).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$checkmark = [char]0x2713

# --- New date header in J4, copying the date-cell formatting from I4 ---
$ws.Range("I4").Copy()
$ws.Range("J4").PasteSpecial(-4122)
$ws.Range("J4").Value = 44335

# --- Fill J5:J15 with a checkmark, copying the per-row formatting from column I ---
foreach ($r in 5..15) {
    $src = $ws.Range("I$r")
    $dst = $ws.Range("J$r")
    $src.Copy()
    $dst.PasteSpecial(-4122)
    $dst.Value = $checkmark
}

# --- Jaime Nunez Delgado (row 9) left early on the 5/12 session ---
$ws.Range("I9").Value = "$checkmark (left at 16:30)"

# --- Column I needs to widen to fit the longer note text ---
$ws.Columns.Item(9).ColumnWidth = 16.333333333333332

# --- Match the author's final selection ---
$ws.Range("J14").Select()
